$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 64
$ws1.Range("F5").Value = 12
$ws1.Range("F6").Value = 161
$ws1.Range("F9").Value = 367
$ws1.Range("F10").Value = 481
$ws1.Range("F13").Value = 12159
$ws1.Range("F14").Value = 5450

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 64
$ws4.Range("F7").Value = 12
$ws4.Range("F8").Value = 161
$ws4.Range("F11").Value = 367
$ws4.Range("F12").Value = 481
$ws4.Range("F15").Value = 12159
$ws4.Range("F17").Value = 5450
